$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet originally holds rows for years 2000, 2002, 2005, 2007, 2010,
# 2012, 2015, 2017 (rows 2-9). The update drops the first four (oldest)
# years and appends a new 2020 row, leaving 2010, 2012, 2015, 2017, 2020
# (rows 2-6).

# Remove the four oldest-year rows (2000年, 2002年, 2005年, 2007年); this
# shifts 2010/2012/2015/2017 up into rows 2-5.
$ws.Range("A2:A5").EntireRow.Delete() | Out-Null

# Copy the formatting of the row above (2017年, now row 5) onto the new
# row 6 so the label cell picks up the same style (bold, bordered,
# centered) as the rest of column A.
$ws.Range("A5").Copy() | Out-Null
$ws.Range("A6").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$excel.CutCopyMode = $false

# Populate the new 2020年 row.
$ws.Range("A6").Value = "2020年"
$ws.Range("B6").Value = 21023911.6478502
$ws.Range("C6").Value = 2377648687.6163
$ws.Range("D6").Value = 441090596.625456
$ws.Range("E6").Value = ""
$ws.Range("F6").Value = 2715857387.03052
$ws.Range("G6").Value = ""
$ws.Range("H6").Value = ""
$ws.Range("I6").Value = 11790811542.2879
$ws.Range("J6").Value = 1918920104.3138
$ws.Range("K6").Value = 263522662.90245
$ws.Range("L6").Value = 118925454.63655
$ws.Range("M6").Value = 314437555.9808
$ws.Range("N6").Value = ""
$ws.Range("O6").Value = 6363715.02579818
$ws.Range("P6").Value = 151082594.831329
$ws.Range("Q6").Value = ""
$ws.Range("R6").Value = 42181555.5967762
$ws.Range("S6").Value = 652874300.461413

Write-Output "done"
